# Rename header columns from *_old / *_new suffixes to *_FV2310 / *_FV2404
# suffixes, convert the data range into a formatted Excel Table (ListObject),
# and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells (row 1) in place: _old -> _FV2310, _new -> _FV2404
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2310")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2404")
        }
    }
}

# Convert the data range A1:U64 into an Excel Table ("Table1") with headers
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1 stays visible while scrolling)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
